# Auto-generated cell updates for Seraph_Profits workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific leve rows
# to refreshed market-board price data, matching the source commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 20 (Leve Item ID / G20 = 1965)
$ws.Range("H20").Value = 760.5
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# ALC row 35 (Leve Item ID / G35 = 1965)
$ws.Range("H35").Value = 760.5
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# ALC row 40 (Leve Item ID / G40 = 5505)
$ws.Range("H40").Value = 1438
$ws.Range("I40").Value = 1378.25
$ws.Range("K40").Value = 1378.25
$ws.Range("M40").Value = -1203.25

# ALC row 100 (Leve Item ID / G100 = 19906)
$ws.Range("H100").Value = 2250
$ws.Range("J100").Value = 2333.3333
$ws.Range("L100").Value = 2333.3333
$ws.Range("N100").Value = -3415.3333

# ALC row 106 (Leve Item ID / G106 = 19903)
$ws.Range("H106").Value = 21473.95
$ws.Range("I106").Value = 22288.053
$ws.Range("J106").Value = 6006
$ws.Range("K106").Value = 22288.053
$ws.Range("L106").Value = 6006
$ws.Range("M106").Value = -21657.053
$ws.Range("N106").Value = -7268

# ALC row 134 (Leve Item ID / G134 = 41997)
$ws.Range("H134").Value = 66666.664
$ws.Range("J134").Value = 66666.664
$ws.Range("L134").Value = 66666.664
$ws.Range("N134").Value = -76806.664

# ALC row 138 (Leve Item ID / G138 = 44169)
$ws.Range("H138").Value = 1971.4482
$ws.Range("I138").Value = 1774.5
$ws.Range("K138").Value = 5323.5
$ws.Range("M138").Value = -183.5

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61 (Leve Item ID / G61 = 43999)
$ws.Range("H61").Value = 5218.7
$ws.Range("I61").Value = 5218.7
$ws.Range("K61").Value = 5218.7
$ws.Range("M61").Value = -5006.7

# ARM row 122 (Leve Item ID / G122 = 36168)
$ws.Range("H122").Value = 785115.25
$ws.Range("I122").Value = 1120731
$ws.Range("K122").Value = 3362193
$ws.Range("M122").Value = -3359743

# ARM row 136 (Leve Item ID / G136 = 43999)
$ws.Range("H136").Value = 5218.7
$ws.Range("I136").Value = 5218.7
$ws.Range("K136").Value = 15656.1
$ws.Range("M136").Value = -13106.1

$ws = $wb.Worksheets.Item("BSM")
# BSM row 11 (Leve Item ID / G11 = 2481)
$ws.Range("H11").Value = 3143.5
$ws.Range("J11").Value = 3449.7144
$ws.Range("L11").Value = 3449.7144
$ws.Range("N11").Value = -3729.7144

# BSM row 30 (Leve Item ID / G30 = 1609)
$ws.Range("H30").Value = 2799.5
$ws.Range("J30").Value = 2799.5
$ws.Range("L30").Value = 2799.5
$ws.Range("N30").Value = -3049.5

# BSM row 42 (Leve Item ID / G42 = 22903)
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58 (Leve Item ID / G58 = 44021)
$ws.Range("H58").Value = 2693.8667
$ws.Range("I58").Value = 1603
$ws.Range("J58").Value = 4875.6
$ws.Range("K58").Value = 1603
$ws.Range("L58").Value = 4875.6
$ws.Range("M58").Value = -1400
$ws.Range("N58").Value = -5281.6

# CRP row 62 (Leve Item ID / G62 = 12580)
$ws.Range("H62").Value = 113788.86
$ws.Range("I62").Value = 66087.164
$ws.Range("K62").Value = 66087.164
$ws.Range("M62").Value = -65463.164

# CRP row 65 (Leve Item ID / G65 = 12580)
$ws.Range("H65").Value = 113788.86
$ws.Range("I65").Value = 66087.164
$ws.Range("K65").Value = 330435.82
$ws.Range("M65").Value = -327315.82

# CRP row 136 (Leve Item ID / G136 = 44021)
$ws.Range("H136").Value = 2693.8667
$ws.Range("I136").Value = 1603
$ws.Range("J136").Value = 4875.6
$ws.Range("K136").Value = 4809
$ws.Range("L136").Value = 14626.8
$ws.Range("M136").Value = -2259
$ws.Range("N136").Value = -19726.8

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4 (Leve Item ID / G4 = 4650)
$ws.Range("H4").Value = 23300342
$ws.Range("I4").Value = 32037064
$ws.Range("J4").Value = 2415.4167
$ws.Range("K4").Value = 96111192
$ws.Range("L4").Value = 7246.250100000001
$ws.Range("M4").Value = -96111080
$ws.Range("N4").Value = -7470.250100000001

# CUL row 11 (Leve Item ID / G11 = 4745)
$ws.Range("H11").Value = 22750
$ws.Range("I11").Value = 22750
$ws.Range("K11").Value = 68250
$ws.Range("M11").Value = -68110

# CUL row 15 (Leve Item ID / G15 = 4661)
$ws.Range("H15").Value = 686.8570999999999
$ws.Range("J15").Value = 867.2857
$ws.Range("L15").Value = 2601.8571
$ws.Range("N15").Value = -2881.8571

# CUL row 26 (Leve Item ID / G26 = 4746)
$ws.Range("H26").Value = 549
$ws.Range("J26").Value = 607.1667
$ws.Range("L26").Value = 1821.5001
$ws.Range("N26").Value = -2397.5001

# CUL row 32 (Leve Item ID / G32 = 4731)
$ws.Range("H32").Value = 6970441.5
$ws.Range("J32").Value = 6970441.5
$ws.Range("L32").Value = 20911324.5
$ws.Range("N32").Value = -20911890.5

# CUL row 39 (Leve Item ID / G39 = 4712)
$ws.Range("H39").Value = 51125
$ws.Range("J39").Value = 100000
$ws.Range("L39").Value = 300000
$ws.Range("N39").Value = -300588

# CUL row 47 (Leve Item ID / G47 = 4663)
$ws.Range("H47").Value = 84.333336
$ws.Range("I47").Value = 24.5
$ws.Range("K47").Value = 73.5
$ws.Range("M47").Value = 357.5

# CUL row 137 (Leve Item ID / G137 = 44088)
$ws.Range("H137").Value = 5116
$ws.Range("J137").Value = 5700
$ws.Range("L137").Value = 17100
$ws.Range("N137").Value = -27300

$ws = $wb.Worksheets.Item("GSM")
# GSM row 31 (Leve Item ID / G31 = 2118)
$ws.Range("H31").Value = 851.4286
$ws.Range("I31").Value = 543.3333
$ws.Range("J31").Value = 2700
$ws.Range("K31").Value = 543.3333
$ws.Range("L31").Value = 2700
$ws.Range("M31").Value = -251.3333
$ws.Range("N31").Value = -3284

# GSM row 37 (Leve Item ID / G37 = 2118)
$ws.Range("H37").Value = 851.4286
$ws.Range("I37").Value = 543.3333
$ws.Range("J37").Value = 2700
$ws.Range("K37").Value = 543.3333
$ws.Range("L37").Value = 2700
$ws.Range("M37").Value = -266.3333
$ws.Range("N37").Value = -3254

# GSM row 132 (Leve Item ID / G132 = 44008)
$ws.Range("H132").Value = 1459.625
$ws.Range("I132").Value = 1239.5714
$ws.Range("K132").Value = 3718.7142
$ws.Range("M132").Value = -1188.7142

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16 (Leve Item ID / G16 = 5289)
$ws.Range("H16").Value = 919.58826
$ws.Range("I16").Value = 945.8125
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 945.8125
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -775.8125
$ws.Range("N16").Value = -840

# LTW row 25 (Leve Item ID / G25 = 3547)
$ws.Range("H25").Value = 29999
$ws.Range("J25").Value = 29999
$ws.Range("L25").Value = 29999
$ws.Range("N25").Value = -30459

# LTW row 46 (Leve Item ID / G46 = 5282)
$ws.Range("H46").Value = 2812.375
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2812.375
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2812.375
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3188.375

# LTW row 122 (Leve Item ID / G122 = 36247)
$ws.Range("H122").Value = 3568
$ws.Range("I122").Value = 2204
$ws.Range("K122").Value = 6612
$ws.Range("M122").Value = -4162

$ws = $wb.Worksheets.Item("WVR")
# WVR row 43 (Leve Item ID / G43 = 3831)
$ws.Range("H43").Value = 11009
$ws.Range("I43").Value = 11009
$ws.Range("K43").Value = 11009
$ws.Range("M43").Value = -10860

# WVR row 49 (Leve Item ID / G49 = 3397)
$ws.Range("H49").Value = 299198.6
$ws.Range("I49").Value = 332665.66
$ws.Range("J49").Value = 248998
$ws.Range("K49").Value = 332665.66
$ws.Range("L49").Value = 248998
$ws.Range("M49").Value = -332435.66
$ws.Range("N49").Value = -249458

# WVR row 81 (Leve Item ID / G81 = 12596)
$ws.Range("H81").Value = 5470.75
$ws.Range("I81").Value = 7737.5
$ws.Range("K81").Value = 15475
$ws.Range("M81").Value = -14414

# WVR row 84 (Leve Item ID / G84 = 12596)
$ws.Range("H84").Value = 5470.75
$ws.Range("I84").Value = 7737.5
$ws.Range("K84").Value = 77375
$ws.Range("M84").Value = -72071

# WVR row 96 (Leve Item ID / G96 = 19977)
$ws.Range("H96").Value = 4816.2
$ws.Range("I96").Value = 4692.4
$ws.Range("J96").Value = 4940
$ws.Range("K96").Value = 4692.4
$ws.Range("L96").Value = 4940
$ws.Range("M96").Value = -3319.4
$ws.Range("N96").Value = -7686

# WVR row 132 (Leve Item ID / G132 = 44029)
$ws.Range("H132").Value = 6750
$ws.Range("I132").Value = 6750
$ws.Range("K132").Value = 20250
$ws.Range("M132").Value = -17720

# WVR row 136 (Leve Item ID / G136 = 44031)
$ws.Range("H136").Value = 980.95
$ws.Range("I136").Value = 728.8333
$ws.Range("K136").Value = 2186.4999
$ws.Range("M136").Value = 363.5001000000002
